# Update cryptos list (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (so numeric-looking
# strings such as "216.60" or "9.02" are not silently coerced into
# Excel numbers), then restore the cell's original style so no stray
# formatting/style index is introduced.
function Set-TextValue {
    param($cell, [string]$value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Straightforward numeric/percentage updates (Price in column D, Volume(1h) in column E)
$updates = @(
    @{ Row = 2;  D = "82.157.75";  E = "  +3.02%  " },
    @{ Row = 3;  D = "3.197.20";   E = "  -0.04%  " },
    @{ Row = 4;  D = $null;        E = "  -0.01%  " },
    @{ Row = 5;  D = "216.60";     E = "  +5.54%  " },
    @{ Row = 6;  D = "625.07";     E = "  -1.57%  " },
    @{ Row = 7;  D = "0.290";      E = "  +22.07%  " },
    @{ Row = 8;  D = $null;        E = "  -0.04%  " },
    @{ Row = 9;  D = $null;        E = "  +0.56%  " },
    @{ Row = 10; D = "3.192.48";   E = "  -0.16%  " },
    @{ Row = 11; D = "0.594";      E = "  +2.51%  " },
    @{ Row = 12; D = $null;        E = "  +11.82%  " },
    @{ Row = 13; D = $null;        E = "  -0.22%  " },
    @{ Row = 14; D = "5.34";       E = "  -3.17%  " },
    @{ Row = 15; D = "3.785.70";   E = "  +0.00%  " },
    @{ Row = 16; D = "31.95";      E = "  +0.28%  " },
    @{ Row = 17; D = "81.934.55";  E = "  +2.99%  " },
    @{ Row = 18; D = "3.191.04";   E = "  -0.39%  " },
    @{ Row = 19; D = "3.25";       E = "  +7.67%  " },
    @{ Row = 20; D = "14.13";      E = "  -2.41%  " },
    @{ Row = 21; D = "436.66";     E = "  +1.69%  " },
    @{ Row = 22; D = "9.02";       E = "  -1.90%  " },
    @{ Row = 23; D = $null;        E = "  +0.61%  " },
    @{ Row = 24; D = "7.27";       E = "  +6.10%  " },
    @{ Row = 25; D = "5.39";       E = "  +13.91%  " },
    @{ Row = 26; D = "3.361.37";   E = "  -0.09%  " },
    @{ Row = 29; D = $null;        E = "  -0.07%  " },
    @{ Row = 30; D = $null;        E = "  +4.33%  " },
    @{ Row = 33; D = "0.998";      E = "  -0.11%  " },
    @{ Row = 34; D = "1.52";       E = "  +1.95%  " },
    @{ Row = 35; D = "0.156";      E = "  +9.76%  " },
    @{ Row = 38; D = "22.88";      E = "  -0.53%  " },
    @{ Row = 39; D = "6.18";       E = "  +11.97%  " },
    @{ Row = 40; D = "0.999";      E = "  -0.02%  " },
    @{ Row = 41; D = $null;        E = "  +0.97%  " },
    @{ Row = 42; D = "2.06";       E = "  +14.45%  " },
    @{ Row = 43; D = "3.09";       E = "  +22.32%  " },
    @{ Row = 44; D = "20.82";      E = "  +3.93%  " },
    @{ Row = 45; D = "161.19";     E = "  -2.51%  " },
    @{ Row = 46; D = $null;        E = "  +0.06%  " },
    @{ Row = 47; D = "188.87";     E = "  -1.80%  " },
    @{ Row = 48; D = $null;        E = "  +3.33%  " },
    @{ Row = 49; D = $null;        E = "  +1.32%  " },
    @{ Row = 50; D = "26.42";      E = "  +2.42%  " },
    @{ Row = 51; D = "0.776";      E = "  -5.57%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextValue $ws.Range("D$r") $u.D
    }
    $ws.Range("E$r").Value = $u.E
}

# Rows that were reordered/swapped (coin moved to a different rank position).
# Row 27 <-> Row 28 : Litecoin / Aptos swap places
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D27") "11.17"
$ws.Range("E27").Value = "  -0.80%  "

$ws.Range("B28").Value = "Litecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D28") "76.74"
$ws.Range("E28").Value = "  -0.45%  "

# Row 31 <-> Row 32 : Bittensor / InternetComputer(DFINITY) swap places
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D31") "9.14"
$ws.Range("E31").Value = "  +1.48%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D32") "587.26"
$ws.Range("E32").Value = "  +11.40%  "

# Row 36 <-> Row 37 : PancakeSwap / Cronos swap places
$ws.Range("B36").Value = "Cronos"
$ws.Range("C36").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D36") "0.140"
$ws.Range("E36").Value = "  +17.31%  "

$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D37") "2.01"
$ws.Range("E37").Value = "  +1.19%  "
